$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.175.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  +2.45%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.870.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.638.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.189.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.42%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.119"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0510"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.308.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.03%  "
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.860"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.12%  "
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.89%  "
$ws.Range("E43").Value = "  -2.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.780.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0964"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.36%  "
